# Fruta / hortaliza, semanal
# Rotate the weekly price-report rows: row 3's data moves to row 4,
# row 4's data moves to row 5, and row 5's data moves up to row 3
# (a cyclic shift of the date/volume/price/origin columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that rotate between rows 3, 4 and 5.
$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $addr3 = "${col}3"
    $addr4 = "${col}4"
    $addr5 = "${col}5"
    $row3[$col] = $ws.Range($addr3).Value2
    $row4[$col] = $ws.Range($addr4).Value2
    $row5[$col] = $ws.Range($addr5).Value2
}

# Apply the cyclic shift: new row3 = old row5, new row4 = old row3, new row5 = old row4
foreach ($col in $cols) {
    $addr3 = "${col}3"
    $addr4 = "${col}4"
    $addr5 = "${col}5"
    $ws.Range($addr3).Value2 = $row5[$col]
    $ws.Range($addr4).Value2 = $row3[$col]
    $ws.Range($addr5).Value2 = $row4[$col]
}
